# Apply cell-level updates to the "cryptos" worksheet to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values are purely numeric-looking strings (e.g. "1.00", "0.800") that
# Excel would otherwise auto-convert to numbers and strip formatting from.
# Force those specific cells to Text format *before* assigning their values so the
# exact original text representation (trailing zeros etc.) is preserved.
$textCells = @("D5", "D6", "D7", "D9", "D15", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D29", "D31", "D32", "D34", "D35", "D38", "D39", "D42", "D43", "D47", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "62.713.29"
$ws.Range("E2").Value = "  +4.47%  "
# Row 3
$ws.Range("D3").Value = "3.343.28"
$ws.Range("E3").Value = "  +4.42%  "
# Row 5
$ws.Range("D5").Value = "561.65"
$ws.Range("E5").Value = "  +4.57%  "
# Row 6
$ws.Range("D6").Value = "152.48"
$ws.Range("E6").Value = "  +4.86%  "
# Row 7
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.13%  "
# Row 8
$ws.Range("E8").Value = "  +0.09%  "
# Row 9
$ws.Range("D9").Value = "7.42"
$ws.Range("E9").Value = "  +0.67%  "
# Row 10
$ws.Range("E10").Value = "  +4.09%  "
# Row 11
$ws.Range("E11").Value = "  +0.26%  "
# Row 12
$ws.Range("D12").Value = "3.919.85"
$ws.Range("E12").Value = "  +4.33%  "
# Row 14
$ws.Range("E14").Value = "  +3.02%  "
# Row 15
$ws.Range("D15").Value = "0.0000180"
$ws.Range("E15").Value = "  +3.42%  "
# Row 16
$ws.Range("D16").Value = "62.718.01"
$ws.Range("E16").Value = "  +4.34%  "
# Row 17
$ws.Range("D17").Value = "3.324.27"
$ws.Range("E17").Value = "  +3.84%  "
# Row 18
$ws.Range("D18").Value = "6.35"
$ws.Range("E18").Value = "  +1.44%  "
# Row 19
$ws.Range("D19").Value = "13.81"
$ws.Range("E19").Value = "  +4.57%  "
# Row 20
$ws.Range("D20").Value = "8.37"
$ws.Range("E20").Value = "  +0.69%  "
# Row 21
$ws.Range("D21").Value = "384.39"
$ws.Range("E21").Value = "  +1.36%  "
# Row 22
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.13%  "
# Row 23
$ws.Range("D23").Value = "0.537"
$ws.Range("E23").Value = "  +2.20%  "
# Row 24
$ws.Range("D24").Value = "70.12"
$ws.Range("E24").Value = "  +0.01%  "
# Row 25
$ws.Range("E25").Value = "  +5.37%  "
# Row 26
$ws.Range("E26").Value = "  -0.18%  "
# Row 27
$ws.Range("D27").Value = "0.0₃0950"
$ws.Range("E27").Value = "  +5.56%  "
# Row 28
$ws.Range("E28").Value = "  +0.04%  "
# Row 29
$ws.Range("D29").Value = "6.55"
$ws.Range("E29").Value = "  +5.63%  "
# Row 30
$ws.Range("E30").Value = "  +3.96%  "
# Row 31
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "5.58"
$ws.Range("E31").Value = "  +2.82%  "
# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "22.91"
$ws.Range("E32").Value = "  +2.51%  "
# Row 33
$ws.Range("E33").Value = "  +7.59%  "
# Row 34
$ws.Range("D34").Value = "6.70"
$ws.Range("E34").Value = "  +0.54%  "
# Row 35
$ws.Range("D35").Value = "160.02"
$ws.Range("E35").Value = "  +1.93%  "
# Row 36
$ws.Range("E36").Value = "  +8.28%  "
# Row 37
$ws.Range("E37").Value = "  +12.70%  "
# Row 38
$ws.Range("D38").Value = "26.94"
$ws.Range("E38").Value = "  +5.29%  "
# Row 39
$ws.Range("D39").Value = "0.0738"
$ws.Range("E39").Value = "  +4.80%  "
# Row 40
$ws.Range("D40").Value = "2.800.28"
$ws.Range("E40").Value = "  +0.10%  "
# Row 41
$ws.Range("E41").Value = "  +6.33%  "
# Row 42
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "0.744"
$ws.Range("E42").Value = "  +3.75%  "
# Row 43
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "40.47"
$ws.Range("E43").Value = "  +1.57%  "
# Row 44
$ws.Range("E44").Value = "  +0.01%  "
# Row 45
$ws.Range("E45").Value = "  +4.15%  "
# Row 46
$ws.Range("D46").Value = "3.387.31"
$ws.Range("E46").Value = "  +4.33%  "
# Row 47
$ws.Range("D47").Value = "21.98"
$ws.Range("E47").Value = "  +6.45%  "
# Row 48
$ws.Range("E48").Value = "  -1.50%  "
# Row 49
$ws.Range("D49").Value = "6.31"
$ws.Range("E49").Value = "  +2.03%  "
# Row 50
$ws.Range("D50").Value = "287.53"
$ws.Range("E50").Value = "  +6.27%  "
# Row 51
$ws.Range("D51").Value = "0.800"
$ws.Range("E51").Value = "  -0.92%  "
